$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) cells that are being updated keep a Text format so that
# numeric-looking strings (e.g. "302.28") are not auto-converted into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.347.52"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "2.325.71"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "302.28"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").Value = "97.95"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "35.64"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "19.60"
$ws.Range("E11").Value = "  +7.34%  "
$ws.Range("D12").Value = "0.0797"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "2.693.53"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "2.343.14"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "43.260.39"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").Value = "12.80"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "68.08"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "237.14"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "25.06"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "2.07"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").Value = "164.21"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "33.23"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("E34").Value = "  +4.45%  "
$ws.Range("D35").Value = "4.50"
$ws.Range("E35").Value = "  -6.90%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "1.988.07"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D43").Value = "19.17"
$ws.Range("E43").Value = "  +7.54%  "
$ws.Range("D44").Value = "10.62"
$ws.Range("E44").Value = "  +5.44%  "
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "2.80"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "2.557.61"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").Value = "54.03"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "72.77"
$ws.Range("E51").Value = "  +0.73%  "
